$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$data = @(
    @("lab.vape.table.setup", "Setup"),
    @("lab.vape.table.mixture", "Mix"),
    @("lab.vape.table.rating", "Hodnocení"),
    @("lab.vape.table.taste", "Chuť"),
    @("lab.vape.table.power", "Výkon"),
    @("lab.vape.table.tc", "Teplota"),
    @("lab.vape.preview.setup", "Setup"),
    @("lab.vape.preview.mixture", "Mix"),
    @("lab.vape.preview.driptip", "Náústek"),
    @("lab.vape.preview.leaks", "Úniky"),
    @("lab.vape.preview.dryhit", "Dryhity"),
    @("lab.vape.preview.rating", "Celkové hodnocení"),
    @("lab.vape.preview.taste", "Hodnocení chuti"),
    @("lab.vape.preview.power", "Výkon (watty)"),
    @("lab.vape.preview.tc", "Teplota"),
    @("lab.vape.preview.airflow", "Airflow"),
    @("lab.vape.preview.juice", "Juice flow"),
    @("lab.vape.preview.mtl", "MTL"),
    @("lab.vape.preview.dl", "DL"),
    @("lab.vape.preview.clouds", "Oblaka"),
    @("lab.vape.preview.fruits", "Ovocné tóny"),
    @("lab.vape.preview.tobacco", "Tabák"),
    @("lab.vape.preview.cakes", "Buchty"),
    @("lab.vape.preview.complex", "Komplexní"),
    @("lab.vape.preview.fresh", "Větrnost"),
    @("lab.vape.preview.atomizer", "Atomizér"),
    @("lab.vape.preview.coil", "Spirálka"),
    @("lab.vape.preview.mod", "Mod")
)

$startRow = 464
$endRow = $startRow + $data.Length - 1

# Copy formatting from the last existing data row (463) onto the new rows first,
# so the new cells share the same style as the rest of the table (style s="1").
$ws.Range("A463:C463").Copy()
$ws.Range("A" + $startRow + ":C" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = "cs"
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

$ws.Cells.Item(485, 2).Select()
$win = $ws.Parent.Windows.Item(1)
$win.ScrollRow = 478
$win.ScrollColumn = 1
